$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2208.5454
$ws.Range("I113").Value = 2132.3333
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 2132.3333
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 1121.6667
$ws.Range("N113").Value = -8808
$ws.Range("H125").Value = 887
$ws.Range("I125").Value = 597.4286
$ws.Range("J125").Value = 1176.5714
$ws.Range("K125").Value = 5376.8574
$ws.Range("L125").Value = 10589.1426
$ws.Range("M125").Value = -2916.8574
$ws.Range("N125").Value = -15509.1426
$ws.Range("H137").Value = 1417.7
$ws.Range("I137").Value = 1063.238
$ws.Range("J137").Value = 2244.7778
$ws.Range("K137").Value = 3189.714
$ws.Range("L137").Value = 6734.3334
$ws.Range("M137").Value = -639.7139999999999
$ws.Range("N137").Value = -11834.3334
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 882
$ws.Range("I2").Value = 721.3
$ws.Range("J2").Value = 1082.875
$ws.Range("K2").Value = 721.3
$ws.Range("L2").Value = 1082.875
$ws.Range("M2").Value = -608.3
$ws.Range("N2").Value = -1308.875
$ws.Range("H32").Value = 847697.1
$ws.Range("I32").Value = 1018090.7
$ws.Range("K32").Value = 1018090.7
$ws.Range("M32").Value = -1017803.7
$ws.Range("H45").Value = 2431
$ws.Range("I45").Value = 1421.1333
$ws.Range("K45").Value = 1421.1333
$ws.Range("M45").Value = -1044.1333
$ws.Range("H63").Value = 4211.316
$ws.Range("I63").Value = 2185.7273
$ws.Range("K63").Value = 2185.7273
$ws.Range("M63").Value = -1499.7273
$ws.Range("H66").Value = 4211.316
$ws.Range("I66").Value = 2185.7273
$ws.Range("K66").Value = 10928.6365
$ws.Range("M66").Value = -7496.636500000001
$ws.Range("H74").Value = 911.9706
$ws.Range("I74").Value = 632.38464
$ws.Range("J74").Value = 1085.0476
$ws.Range("K74").Value = 632.38464
$ws.Range("L74").Value = 1085.0476
$ws.Range("M74").Value = 241.61536
$ws.Range("N74").Value = -2833.0476
$ws.Range("H77").Value = 911.9706
$ws.Range("I77").Value = 632.38464
$ws.Range("J77").Value = 1085.0476
$ws.Range("K77").Value = 3161.9232
$ws.Range("L77").Value = 5425.238
$ws.Range("M77").Value = 1206.0768
$ws.Range("N77").Value = -14161.238
$ws.Range("H116").Value = 882
$ws.Range("I116").Value = 721.3
$ws.Range("J116").Value = 1082.875
$ws.Range("K116").Value = 721.3
$ws.Range("L116").Value = 1082.875
$ws.Range("M116").Value = 1572.7
$ws.Range("N116").Value = -5670.875
$ws.Range("H122").Value = 1497.5385
$ws.Range("I122").Value = 1447.8
$ws.Range("J122").Value = 1663.3334
$ws.Range("K122").Value = 4343.4
$ws.Range("L122").Value = 4990.0002
$ws.Range("M122").Value = -1893.4
$ws.Range("N122").Value = -9890.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 882
$ws.Range("I3").Value = 721.3
$ws.Range("J3").Value = 1082.875
$ws.Range("K3").Value = 721.3
$ws.Range("L3").Value = 1082.875
$ws.Range("M3").Value = -607.3
$ws.Range("N3").Value = -1310.875
$ws.Range("H86").Value = 1421.5238
$ws.Range("I86").Value = 1310.4
$ws.Range("J86").Value = 1699.3334
$ws.Range("K86").Value = 1310.4
$ws.Range("L86").Value = 1699.3334
$ws.Range("M86").Value = -187.4000000000001
$ws.Range("N86").Value = -3945.3334
$ws.Range("H89").Value = 1421.5238
$ws.Range("I89").Value = 1310.4
$ws.Range("J89").Value = 1699.3334
$ws.Range("K89").Value = 6552
$ws.Range("L89").Value = 8496.666999999999
$ws.Range("M89").Value = -936
$ws.Range("N89").Value = -19728.667
$ws.Range("H99").Value = 1028.0454
$ws.Range("I99").Value = 950.94446
$ws.Range("J99").Value = 1375
$ws.Range("K99").Value = 950.94446
$ws.Range("L99").Value = 1375
$ws.Range("M99").Value = 547.05554
$ws.Range("N99").Value = -4371
$ws.Range("H107").Value = 1988.4546
$ws.Range("I107").Value = 1874.7778
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 1874.7778
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 45.22219999999993
$ws.Range("N107").Value = -6340
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4173.451
$ws.Range("I31").Value = 1032.48
$ws.Range("J31").Value = 7193.615
$ws.Range("K31").Value = 1032.48
$ws.Range("L31").Value = 7193.615
$ws.Range("M31").Value = -737.48
$ws.Range("N31").Value = -7783.615
$ws.Range("H34").Value = 4173.451
$ws.Range("I34").Value = 1032.48
$ws.Range("J34").Value = 7193.615
$ws.Range("K34").Value = 1032.48
$ws.Range("L34").Value = 7193.615
$ws.Range("M34").Value = -830.48
$ws.Range("N34").Value = -7597.615
$ws.Range("H107").Value = 2718153.5
$ws.Range("I107").Value = 6250477
$ws.Range("J107").Value = 981.46155
$ws.Range("K107").Value = 6250477
$ws.Range("L107").Value = 981.46155
$ws.Range("M107").Value = -6248557
$ws.Range("N107").Value = -4821.46155
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 771.26666
$ws.Range("I5").Value = 400.6279
$ws.Range("J5").Value = 1269.3125
$ws.Range("K5").Value = 1201.8837
$ws.Range("L5").Value = 3807.9375
$ws.Range("M5").Value = -1089.8837
$ws.Range("N5").Value = -4031.9375
$ws.Range("H64").Value = 1350.5
$ws.Range("J64").Value = 2001
$ws.Range("L64").Value = 6003
$ws.Range("N64").Value = -6543
$ws.Range("H67").Value = 1350.5
$ws.Range("J67").Value = 2001
$ws.Range("L67").Value = 6003
$ws.Range("N67").Value = -7875
$ws.Range("H68").Value = 1115.9531
$ws.Range("J68").Value = 1178.8334
$ws.Range("L68").Value = 3536.5002
$ws.Range("N68").Value = -5158.5002
$ws.Range("H71").Value = 1115.9531
$ws.Range("J71").Value = 1178.8334
$ws.Range("L71").Value = 10609.5006
$ws.Range("N71").Value = -18721.5006
$ws.Range("H110").Value = 10741.777
$ws.Range("I110").Value = 1919
$ws.Range("J110").Value = 17800
$ws.Range("K110").Value = 5757
$ws.Range("L110").Value = 53400
$ws.Range("M110").Value = -1667
$ws.Range("N110").Value = -61580
$ws.Range("H131").Value = 1139.5869
$ws.Range("I131").Value = 1025.8462
$ws.Range("J131").Value = 1184.3939
$ws.Range("K131").Value = 3077.5386
$ws.Range("L131").Value = 3553.1817
$ws.Range("M131").Value = 1962.4614
$ws.Range("N131").Value = -13633.1817
$ws.Range("H135").Value = 771.26666
$ws.Range("I135").Value = 400.6279
$ws.Range("J135").Value = 1269.3125
$ws.Range("K135").Value = 3605.6511
$ws.Range("L135").Value = 11423.8125
$ws.Range("M135").Value = -1070.6511
$ws.Range("N135").Value = -16493.8125
$ws.Range("H136").Value = 3222
$ws.Range("I136").Value = 3010.8
$ws.Range("K136").Value = 9032.400000000001
$ws.Range("M136").Value = -3932.400000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2320.6667
$ws.Range("I132").Value = 978.5714
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 2935.7142
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -405.7142000000003
$ws.Range("N132").Value = -17658.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7985.8
$ws.Range("I22").Value = 660
$ws.Range("J22").Value = 16358.143
$ws.Range("K22").Value = 660
$ws.Range("L22").Value = 16358.143
$ws.Range("M22").Value = -365
$ws.Range("N22").Value = -16948.143
$ws.Range("H27").Value = 7985.8
$ws.Range("I27").Value = 660
$ws.Range("J27").Value = 16358.143
$ws.Range("K27").Value = 660
$ws.Range("L27").Value = 16358.143
$ws.Range("M27").Value = -553
$ws.Range("N27").Value = -16572.143
$ws.Range("H68").Value = 1362.3077
$ws.Range("J68").Value = 1600
$ws.Range("L68").Value = 1600
$ws.Range("N68").Value = -3098
$ws.Range("H71").Value = 1362.3077
$ws.Range("J71").Value = 1600
$ws.Range("L71").Value = 8000
$ws.Range("N71").Value = -15488
$ws.Range("H93").Value = 7224.294
$ws.Range("I93").Value = 9058.75
$ws.Range("J93").Value = 2821.6
$ws.Range("K93").Value = 9058.75
$ws.Range("L93").Value = 2821.6
$ws.Range("M93").Value = -7810.75
$ws.Range("N93").Value = -5317.6
$ws.Range("H100").Value = 2843.1428
$ws.Range("I100").Value = 2814
$ws.Range("J100").Value = 2859.3333
$ws.Range("K100").Value = 2814
$ws.Range("L100").Value = 2859.3333
$ws.Range("M100").Value = -2273
$ws.Range("N100").Value = -3941.3333
$ws.Range("H122").Value = 3366.6667
$ws.Range("I122").Value = 3216.6667
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 9650.000100000001
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -7200.000100000001
$ws.Range("N122").Value = -15900.0001
$ws.Range("H132").Value = 3540.712
$ws.Range("I132").Value = 3787
$ws.Range("J132").Value = 3332.9062
$ws.Range("K132").Value = 11361
$ws.Range("L132").Value = 9998.7186
$ws.Range("M132").Value = -8831
$ws.Range("N132").Value = -15058.7186
$ws.Range("H136").Value = 1363.1111
$ws.Range("I136").Value = 1269.2
$ws.Range("J136").Value = 1631.4286
$ws.Range("K136").Value = 3807.6
$ws.Range("L136").Value = 4894.2858
$ws.Range("M136").Value = -1257.6
$ws.Range("N136").Value = -9994.2858
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8336308
$ws.Range("I132").Value = 3728.353
$ws.Range("J132").Value = 16205966
$ws.Range("K132").Value = 11185.059
$ws.Range("L132").Value = 48617898
$ws.Range("M132").Value = -8655.059000000001
$ws.Range("N132").Value = -48622958
